# Update cryptos list - refresh Price (D) and Volume(1h) (E) columns, and
# swap the Cronos / VeChain rows (43 and 44) to reflect the new ranking.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => (Price, Volume) ; $null means "leave value unchanged"
$updates = @{
    2  = @("36.659.95",  "  -2.06%  ")
    3  = @("1.988.91",   "  -2.32%  ")
    4  = @($null,        "  +0.00%  ")
    5  = @("256.13",     "  +3.33%  ")
    6  = @("0.611",      "  -2.35%  ")
    7  = @($null,        "  +0.16%  ")
    8  = @("55.32",      "  -8.05%  ")
    9  = @("0.375",      "  -4.94%  ")
    10 = @("0.0762",     "  -5.88%  ")
    11 = @("0.101",      "  -3.07%  ")
    12 = @("14.26",      "  -6.23%  ")
    13 = @("2.284.50",   "  -2.13%  ")
    14 = @("21.24",      "  -4.31%  ")
    15 = @("0.788",      "  -7.84%  ")
    16 = @("5.15",       "  -5.85%  ")
    17 = @("1.998.03",   "  -1.80%  ")
    18 = @("36.541.68",  "  -2.11%  ")
    19 = @("70.23",      "  -0.45%  ")
    20 = @("0.0₃0820",   "  -4.99%  ")
    21 = @("234.66",     "  +1.70%  ")
    22 = @("5.04",       "  -4.15%  ")
    23 = @($null,        "  -0.25%  ")
    24 = @("2.53",       "  -1.03%  ")
    25 = @("2.36",       "  +0.33%  ")
    26 = @("163.80",     "  -0.38%  ")
    27 = @("8.83",       "  -6.21%  ")
    28 = @("19.30",      "  -3.53%  ")
    29 = @("1.34",       "  -3.20%  ")
    30 = @($null,        "  -10.01%  ")
    31 = @("0.118",      "  -3.01%  ")
    32 = @("4.52",       "  -5.71%  ")
    33 = @($null,        "  -7.31%  ")
    34 = @($null,        "  -3.66%  ")
    35 = @("2.35",       "  -8.68%  ")
    36 = @("3.43",       "  -5.71%  ")
    37 = @("1.81",       "  +0.52%  ")
    38 = @($null,        "  -0.01%  ")
    39 = @("5.48",       "  +0.30%  ")
    40 = @($null,        "  +0.54%  ")
    41 = @("1.441.73",   "  +4.30%  ")
    42 = @("1.16",       "  -1.59%  ")
    45 = @("88.43",      "  -3.76%  ")
    46 = @("15.51",      "  -7.63%  ")
    47 = @($null,        "  -4.55%  ")
    48 = @($null,        "  -0.29%  ")
    49 = @("6.87",       "  -9.55%  ")
    50 = @("2.176.63",   "  -2.14%  ")
    51 = @($null,        "  -9.21%  ")
}

foreach ($row in $updates.Keys) {
    $pair = $updates[$row]
    $price = $pair[0]
    $volume = $pair[1]
    if ($null -ne $price) {
        $ws.Range("D$row").Value = $price
    }
    $ws.Range("E$row").Value = $volume
}

# Rows 43/44 swap: Cronos moves above VeChain in the ranking.
$ws.Range("B43").Value = "Cronos"
$ws.Range("C43").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D43").Value = "0.0914"
$ws.Range("E43").Value = "  -6.16%  "

$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "0.0207"
$ws.Range("E44").Value = "  -4.87%  "
